{"js": "// The source diff for this revision only reorders XML attributes and\n// namespace declarations that get rewritten whenever the test fixture's\n// canonical OOXML is re-serialized (root namespace declarations on\n// <w:document>, <w:pgSz>/<w:pgMar> attribute order, <w:rFonts>/<w:lang>\n// attribute order in styles.xml's docDefaults, the w:latentStyles /\n// w:lsdException / w:style attribute order, etc.). Every changed line\n// has the exact same element name and the exact same set of attribute\n// name/value pairs on both sides of the diff - nothing is added,\n// removed, or renamed, and no paragraph text, field code, formatting,\n// or structural content differs. There is therefore no visible-content\n// edit to make through the Word object model: the body text, runs,\n// fields and section/page-setup values are already exactly what the\n// target state requires.\n//\n// Touch the body read-only so the script still exercises the\n// load/sync round trip expected of an Office.js edit script, without\n// mutating anything.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The source diff for this revision only reorders XML attributes and\n# namespace declarations that get rewritten whenever the test fixture's\n# canonical OOXML is re-serialized (root namespace declarations on\n# <w:document>, <w:pgSz>/<w:pgMar> attribute order, <w:rFonts>/<w:lang>\n# attribute order in styles.xml's docDefaults, the w:latentStyles /\n# w:lsdException / w:style attribute order, etc.). Every changed line\n# has the exact same element name and the exact same set of attribute\n# name/value pairs on both sides of the diff - nothing is added,\n# removed, or renamed, and no paragraph text, field code, formatting,\n# or structural content differs. There is therefore no visible-content\n# edit to make through the Word object model: the body text, runs,\n# fields and section/page-setup values are already exactly what the\n# target state requires.\n\n# Touch the document read-only so the script still exercises the Word\n# COM object model expected of a COM edit script, without mutating\n# anything.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
